# Updates cryptos list values (Price / Volume(1h) columns) per the
# upstream GitHub Actions scrape-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text so Excel does not silently
# coerce numeric-looking strings (e.g. '0.340', '1.00') into numbers
# and strip trailing zeros, while keeping the cell's original style
# (no NumberFormat residue left behind).
function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '58.828.26'
Set-TextValue "E2" '  +0.44%  '
Set-TextValue "D3" '2.572.52'
Set-TextValue "E3" '  -0.36%  '
Set-TextValue "E4" '  +0.06%  '
Set-TextValue "D5" '562.03'
Set-TextValue "E5" '  +3.55%  '
Set-TextValue "D6" '142.47'
Set-TextValue "E6" '  -1.14%  '
Set-TextValue "E7" '  +0.07%  '
Set-TextValue "E8" '  +2.15%  '
Set-TextValue "D9" '2.577.86'
Set-TextValue "E9" '  -0.43%  '
Set-TextValue "E10" '  -2.00%  '
Set-TextValue "E11" '  +2.33%  '
Set-TextValue "E12" '  +8.52%  '
Set-TextValue "D13" '0.340'
Set-TextValue "E13" '  +2.21%  '
Set-TextValue "D14" '3.024.13'
Set-TextValue "E14" '  -0.35%  '
Set-TextValue "D15" '58.933.02'
Set-TextValue "E15" '  +0.76%  '
Set-TextValue "D16" '21.87'
Set-TextValue "E16" '  +6.11%  '
Set-TextValue "E17" '  +3.50%  '
Set-TextValue "D18" '2.587.59'
Set-TextValue "E18" '  +0.13%  '
Set-TextValue "D19" '4.49'
Set-TextValue "E19" '  +0.83%  '
Set-TextValue "D20" '334.67'
Set-TextValue "E20" '  +0.21%  '
Set-TextValue "D21" '10.14'
Set-TextValue "E21" '  +0.90%  '
Set-TextValue "D22" '6.14'
Set-TextValue "E22" '  +0.90%  '
Set-TextValue "D23" '0.999'
Set-TextValue "E23" '  -0.10%  '
Set-TextValue "D24" '64.49'
Set-TextValue "E24" '  -2.81%  '
Set-TextValue "E25" '  +4.85%  '
Set-TextValue "D26" '1.00'
Set-TextValue "E26" '  +0.41%  '
Set-TextValue "D27" '0.160'
Set-TextValue "E27" '  +1.68%  '
Set-TextValue "D28" '7.20'
Set-TextValue "E28" '  +1.39%  '
Set-TextValue "D29" '0.0₃0777'
Set-TextValue "E29" '  +4.89%  '
Set-TextValue "E30" '  +0.02%  '
Set-TextValue "E31" '  +1.96%  '
Set-TextValue "D32" '160.63'
Set-TextValue "E32" '  +5.15%  '
Set-TextValue "E33" '  +1.26%  '
Set-TextValue "D34" '18.87'
Set-TextValue "E34" '  -0.23%  '
Set-TextValue "E35" '  +2.15%  '
Set-TextValue "D36" '0.876'
Set-TextValue "E36" '  +3.09%  '
Set-TextValue "D37" '0.876'
Set-TextValue "E37" '  +6.44%  '
Set-TextValue "E38" '  +2.96%  '
Set-TextValue "D39" '36.69'
Set-TextValue "E39" '  -1.09%  '
Set-TextValue "D40" '1.47'
Set-TextValue "E40" '  +3.73%  '
Set-TextValue "D41" '295.23'
Set-TextValue "E42" '  +0.66%  '
Set-TextValue "E43" '  +0.06%  '
Set-TextValue "D44" '0.0969'
Set-TextValue "E44" '  +2.97%  '
Set-TextValue "D45" '0.593'
Set-TextValue "E45" '  +0.11%  '
Set-TextValue "E46" '  +1.30%  '
Set-TextValue "E47" '  -0.23%  '
Set-TextValue "D48" '125.30'
Set-TextValue "E48" '  +14.55%  '
Set-TextValue "D49" '18.96'
Set-TextValue "E49" '  +2.34%  '
Set-TextValue "E50" '  +1.71%  '
Set-TextValue "D51" '18.39'
Set-TextValue "E51" '  +2.88%  '
